$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove the training-metric columns (E:H); this shifts test -> E:H, vs -> I:L, Kernel -> M
$ws.Range("E:H").Delete()

# Step 2: write back the recalculated metric values and kernel labels exactly as produced by the corrected run
$ws.Range("A2").Value = 0.7343635235672223
$ws.Range("B2").Value = 0.1576932873087903
$ws.Range("C2").Value = 0.2322280469286339
$ws.Range("D2").Value = 0.3971061411119076
$ws.Range("E2").Value = 0.902823813176789
$ws.Range("F2").Value = 0.08276894642775777
$ws.Range("G2").Value = 0.2156540326768729
$ws.Range("H2").Value = 0.2876959270267095
$ws.Range("I2").Value = 0.7409829035085965
$ws.Range("J2").Value = 0.1618975137598348
$ws.Range("K2").Value = 0.1896063190244313
$ws.Range("L2").Value = 0.4023649012523767
$ws.Range("M2").Value = "RBF"

$ws.Range("A3").Value = 0.7540483215612064
$ws.Range("B3").Value = 0.1460075408805644
$ws.Range("C3").Value = 0.248954857745756
$ws.Range("D3").Value = 0.3821093310566551
$ws.Range("E3").Value = 0.8783149841011321
$ws.Range("F3").Value = 0.1036441219937698
$ws.Range("G3").Value = 0.2790110877488802
$ws.Range("H3").Value = 0.321938071674926
$ws.Range("I3").Value = 0.7540739105741914
$ws.Range("J3").Value = 0.1537150365209141
$ws.Range("K3").Value = 0.2197416489671756
$ws.Range("L3").Value = 0.3920650921988773
$ws.Range("M3").Value = "Matern_0.5"

$ws.Range("A4").Value = 0.7343638965567995
$ws.Range("B4").Value = 0.157693065886065
$ws.Range("C4").Value = 0.2322279869070913
$ws.Range("D4").Value = 0.3971058623164168
$ws.Range("E4").Value = 0.9028239517586609
$ws.Range("F4").Value = 0.08276882839188993
$ws.Range("G4").Value = 0.215653964950414
$ws.Range("H4").Value = 0.2876957218866661
$ws.Range("I4").Value = 0.740983340364318
$ws.Range("J4").Value = 0.1618972407050592
$ws.Range("K4").Value = 0.1896061687584276
$ws.Range("L4").Value = 0.4023645619398647
$ws.Range("M4").Value = "RationalQuadratic"

$ws.Range("A5").Value = 0.7343633781643774
$ws.Range("B5").Value = 0.1576933736262001
$ws.Range("C5").Value = 0.2322280673868489
$ws.Range("D5").Value = 0.3971062497949385
$ws.Range("E5").Value = 0.9028237842013376
$ws.Range("F5").Value = 0.08276897110733886
$ws.Range("G5").Value = 0.2156540925739203
$ws.Range("H5").Value = 0.2876959699184868
$ws.Range("I5").Value = 0.7409827482907532
$ws.Range("J5").Value = 0.1618976107780753
$ws.Range("K5").Value = 0.1896061378365506
$ws.Range("L5").Value = 0.4023650218123778
$ws.Range("M5").Value = "ExpSineSquared"

$ws.Range("A6").Value = -0.03181051529939483
$ws.Range("B6").Value = 0.61252729377516
$ws.Range("C6").Value = 0.5797516786087576
$ws.Range("D6").Value = 0.7826412292839932
$ws.Range("E6").Value = 0.00001123448542794669
$ws.Range("F6").Value = 0.8517314711248334
$ws.Range("G6").Value = 0.8541377576182823
$ws.Range("H6").Value = 0.9228929900724316
$ws.Range("I6").Value = -0.0106886765959433
$ws.Range("J6").Value = 0.6317265776760479
$ws.Range("K6").Value = 0.5357825163404412
$ws.Range("L6").Value = 0.7948122908435978
$ws.Range("M6").Value = "DotProduct"

